$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the row for student 232004 (row 3) so the row for 232005 shifts up
$ws.Rows.Item(3).Delete()

# Update the "User" column (F) values for the remaining data rows
$ws.Range("F2").Value = "admin@admin.com"
$ws.Range("F3").Value = "admin@admin.com"

$wb.Save()
